$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 151 from placeholder text to numeric data, and add rows 152-265
$values = @(
  @(151, 27, 67),
  @(152, 27, 67),
  @(153, 27, 67),
  @(154, 27, 66),
  @(155, 27, 66),
  @(156, 27, 66),
  @(157, 27, 66),
  @(158, 27, 66),
  @(159, 27, 66),
  @(160, 27, 65),
  @(161, 27, 65),
  @(162, 27, 65),
  @(163, 27, 64),
  @(164, 27, 64),
  @(165, 27, 63),
  @(166, 27, 65),
  @(167, 27, 64),
  @(168, 27, 64),
  @(169, 27, 63),
  @(170, 27, 63),
  @(171, 27, 63),
  @(172, 27, 63),
  @(173, 27, 63),
  @(174, 27, 62),
  @(175, 27, 62),
  @(176, 27, 62),
  @(177, 25, 43),
  @(178, 25, 43),
  @(179, 25, 43),
  @(180, 25, 43),
  @(181, 25, 43),
  @(182, 25, 43),
  @(183, 25, 43),
  @(184, 25, 42),
  @(185, 25, 42),
  @(186, 25, 42),
  @(187, 25, 43),
  @(188, 25, 43),
  @(189, 25, 43),
  @(190, 25, 43),
  @(191, 25, 43),
  @(192, 25, 43),
  @(193, 25, 43),
  @(194, 25, 43),
  @(195, 25, 43),
  @(196, 25, 42),
  @(197, 25, 43),
  @(198, 25, 43),
  @(199, 25, 43),
  @(200, 25, 43),
  @(201, 25, 43),
  @(202, 25, 42),
  @(203, 25, 42),
  @(204, 25, 42),
  @(205, 25, 42),
  @(206, 25, 42),
  @(207, 25, 42),
  @(208, 25, 42),
  @(209, 25, 42),
  @(210, 25, 42),
  @(211, 25, 41),
  @(212, 25, 41),
  @(213, 25, 41),
  @(214, 25, 41),
  @(215, 25, 41),
  @(216, 25, 41),
  @(217, 25, 41),
  @(218, 25, 41),
  @(219, 25, 41),
  @(220, 25, 41),
  @(221, 25, 41),
  @(222, 25, 41),
  @(223, 25, 41),
  @(224, 25, 41),
  @(225, 25, 40),
  @(226, 25, 40),
  @(227, 25, 40),
  @(228, 25, 40),
  @(229, 25, 40),
  @(230, 25, 40),
  @(231, 25, 39),
  @(232, 25, 40),
  @(233, 25, 40),
  @(234, 25, 40),
  @(235, 25, 40),
  @(236, 25, 41),
  @(237, 25, 40),
  @(238, 25, 40),
  @(239, 25, 40),
  @(240, 25, 40),
  @(241, 25, 40),
  @(242, 25, 40),
  @(243, 25, 40),
  @(244, 25, 40),
  @(245, 25, 40),
  @(246, 25, 40),
  @(247, 25, 40),
  @(248, 25, 40),
  @(249, 25, 40),
  @(250, 25, 39),
  @(251, 25, 39),
  @(252, 25, 39),
  @(253, 25, 39),
  @(254, 25, 39),
  @(255, 25, 39),
  @(256, 25, 39),
  @(257, 25, 39),
  @(258, 25, 39),
  @(259, 25, 39),
  @(260, 25, 38),
  @(261, 25, 38),
  @(262, 25, 38),
  @(263, 25, 38),
  @(264, 25, 37),
  @(265, 25, 36)
)

foreach ($row in $values) {
  $r = $row[0]
  $a = $row[1]
  $b = $row[2]
  $ws.Cells.Item($r, 1).Value = $a
  $ws.Cells.Item($r, 2).Value = $b
}

# Move the trailing placeholder (blank-ish) text row to the new last row (266)
$ws.Cells.Item(266, 1).Value = " "
$ws.Cells.Item(266, 2).Value = "  "

# Update the view: scroll position and active selection
$ws.Range("C177").Select()
